# Apply changes described by the diff:
# 1. Rename sheet "SCD0205" -> "SCD0011"
# 2. Update B2 (TC_ID) value from "DGS-220" to "SCD0011-036"
# 3. Widen column B to fit new text
# 4. Move the active selection to B3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "SCD0011"

# Update the TC_ID cell value
$ws.Range("B2").Value = "SCD0011-036"

# Adjust column B width to fit the new (wider) text
$ws.Columns("B").ColumnWidth = 11.712

# Update the active cell selection
$ws.Range("B3").Select()
